$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4,D5,D6,D7,D12,D14,D16,D20,D21,D22,D23,D24,D28,D29,D30,D32,D33,D34,D37,D38,D40,D41,D43,D45,D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.408.93"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.288.93"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "576.81"
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").Value = "182.56"
$ws.Range("E6").Value = "  -2.93%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.285.26"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("E11").Value = "  -2.63%  "
$ws.Range("D12").Value = "46.41"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "633.05"
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("D15").Value = "3.812.39"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "8.39"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "65.576.56"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "3.287.98"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "17.61"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "0.885"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "17.92"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "99.98"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "9.33"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "30.58"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "8.33"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "573.71"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").Value = "3.69"
$ws.Range("E33").Value = "  -9.05%  "
$ws.Range("D34").Value = "10.84"
$ws.Range("E34").Value = "  -1.95%  "
$ws.Range("D35").Value = "3.841.11"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "55.45"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "32.34"
$ws.Range("E40").Value = "  -4.57%  "
$ws.Range("B41").Value = "ApeXProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D41").Value = "3.40"
$ws.Range("E41").Value = "  +4.88%  "
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.58"
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0677"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("E50").Value = "  -2.85%  "
$ws.Range("D51").Value = "129.24"
$ws.Range("E51").Value = "  +5.45%  "
